$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.5110453333333334
$ws.Range("H2").Value = 1.533136
$ws.Range("I2").Value = 0.1569529625135799
$ws.Range("J2").Value = 0.1569529625135799
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 12.28101533333333
$ws.Range("N2").Value = 36.843046
$ws.Range("O2").Value = 0.959552102275422
$ws.Range("P2").Value = 0.959552102275422
$ws.Range("Q2").Value = 6.276155574695112
$ws.Range("R2").Value = 56.48540017225601
$ws.Range("S2").Value = 0.1506045451382611
$ws.Range("T2").Value = 0.1506045451382611

$ws.Range("G3").Value = 0.5110453333333334
$ws.Range("H3").Value = 1.533136
$ws.Range("I3").Value = 0.1569529625135799
$ws.Range("J3").Value = 0.1569529625135799
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.09168666666666665
$ws.Range("N3").Value = 0.27506
$ws.Range("O3").Value = 0.007163750827004844
$ws.Range("P3").Value = 0.007163750827004845
$ws.Range("Q3").Value = 0.04685604312888889
$ws.Range("R3").Value = 0.4217043881599999
$ws.Range("S3").Value = 0.001124371915007518
$ws.Range("T3").Value = 0.001124371915007518

$ws.Range("G4").Value = 0.5110453333333334
$ws.Range("H4").Value = 1.533136
$ws.Range("I4").Value = 0.1569529625135799
$ws.Range("J4").Value = 0.1569529625135799
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4259936666666666
$ws.Range("N4").Value = 1.277981
$ws.Range("O4").Value = 0.03328414689757318
$ws.Range("P4").Value = 0.03328414689757318
$ws.Range("Q4").Value = 0.2177020753795555
$ws.Range("R4").Value = 1.959318678416
$ws.Range("S4").Value = 0.005224045460311288
$ws.Range("T4").Value = 0.005224045460311288

$ws.Range("G5").Value = 1.884975666666667
$ws.Range("H5").Value = 5.654927
$ws.Range("I5").Value = 0.5789163814873767
$ws.Range("J5").Value = 0.5789163814873767
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 12.28101533333333
$ws.Range("N5").Value = 36.843046
$ws.Range("O5").Value = 0.959552102275422
$ws.Range("P5").Value = 0.959552102275422
$ws.Range("Q5").Value = 23.14941506529356
$ws.Range("R5").Value = 208.344735587642
$ws.Range("S5").Value = 0.5555004308978925
$ws.Range("T5").Value = 0.5555004308978925

$ws.Range("G6").Value = 1.884975666666667
$ws.Range("H6").Value = 5.654927
$ws.Range("I6").Value = 0.5789163814873767
$ws.Range("J6").Value = 0.5789163814873767
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.09168666666666665
$ws.Range("N6").Value = 0.27506
$ws.Range("O6").Value = 0.007163750827004844
$ws.Range("P6").Value = 0.007163750827004845
$ws.Range("Q6").Value = 0.1728271356244444
$ws.Range("R6").Value = 1.55544422062
$ws.Range("S6").Value = 0.004147212706646846
$ws.Range("T6").Value = 0.004147212706646847

$ws.Range("G7").Value = 1.884975666666667
$ws.Range("H7").Value = 5.654927
$ws.Range("I7").Value = 0.5789163814873767
$ws.Range("J7").Value = 0.5789163814873767
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.4259936666666666
$ws.Range("N7").Value = 1.277981
$ws.Range("O7").Value = 0.03328414689757318
$ws.Range("P7").Value = 0.03328414689757318
$ws.Range("Q7").Value = 0.8029876958207777
$ws.Range("R7").Value = 7.226889262386999
$ws.Range("S7").Value = 0.01926873788283736
$ws.Range("T7").Value = 0.01926873788283736

$ws.Range("G8").Value = 0.8600203333333334
$ws.Range("H8").Value = 2.580061
$ws.Range("I8").Value = 0.2641306559990434
$ws.Range("J8").Value = 0.2641306559990435
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 12.28101533333333
$ws.Range("N8").Value = 36.843046
$ws.Range("O8").Value = 0.959552102275422
$ws.Range("P8").Value = 0.959552102275422
$ws.Range("Q8").Value = 10.56192290064511
$ws.Range("R8").Value = 95.05730610580601
$ws.Range("S8").Value = 0.2534471262392684
$ws.Range("T8").Value = 0.2534471262392685

$ws.Range("G9").Value = 0.8600203333333334
$ws.Range("H9").Value = 2.580061
$ws.Range("I9").Value = 0.2641306559990434
$ws.Range("J9").Value = 0.2641306559990435
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.09168666666666665
$ws.Range("N9").Value = 0.27506
$ws.Range("O9").Value = 0.007163750827004844
$ws.Range("P9").Value = 0.007163750827004845
$ws.Range("Q9").Value = 0.07885239762888888
$ws.Range("R9").Value = 0.70967157866
$ws.Range("S9").Value = 0.001892166205350479
$ws.Range("T9").Value = 0.00189216620535048

$ws.Range("G10").Value = 0.8600203333333334
$ws.Range("H10").Value = 2.580061
$ws.Range("I10").Value = 0.2641306559990434
$ws.Range("J10").Value = 0.2641306559990435
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.4259936666666666
$ws.Range("N10").Value = 1.277981
$ws.Range("O10").Value = 0.03328414689757318
$ws.Range("P10").Value = 0.03328414689757318
$ws.Range("Q10").Value = 0.3663632152045555
$ws.Range("R10").Value = 3.297268936841
$ws.Range("S10").Value = 0.008791363554424529
$ws.Range("T10").Value = 0.008791363554424531
